# Apply F-column numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 20301
$ws.Range("F5").Value = 321
$ws.Range("F6").Value = 1106
$ws.Range("F8").Value = 7644
$ws.Range("F9").Value = 523
$ws.Range("F10").Value = 738
$ws.Range("F11").Value = 276
$ws.Range("F12").Value = 42
$ws.Range("F13").Value = 163
$ws.Range("F14").Value = 130
$ws.Range("F17").Value = 199
$ws.Range("F19").Value = 450
$ws.Range("F20").Value = 75
$ws.Range("F21").Value = 689
$ws.Range("F24").Value = 71
$ws.Range("F25").Value = 327
$ws.Range("F26").Value = 1122
$ws.Range("F27").Value = 36
$ws.Range("F28").Value = 20
$ws.Range("F29").Value = 187
$ws.Range("F31").Value = 570
$ws.Range("F32").Value = 82
$ws.Range("F33").Value = 2878
$ws.Range("F34").Value = 26
$ws.Range("F37").Value = 12713
$ws.Range("F38").Value = 1341
$ws.Range("F40").Value = 32
$ws.Range("F41").Value = 61
$ws.Range("F42").Value = 276
$ws.Range("F43").Value = 381
$ws.Range("F44").Value = 4012
$ws.Range("F46").Value = 96

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 203

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 20301
$ws.Range("F5").Value = 321
$ws.Range("F6").Value = 1106
$ws.Range("F8").Value = 7644
$ws.Range("F9").Value = 523
$ws.Range("F10").Value = 738
$ws.Range("F11").Value = 276
$ws.Range("F12").Value = 42
$ws.Range("F13").Value = 163
$ws.Range("F14").Value = 130
$ws.Range("F17").Value = 199
$ws.Range("F19").Value = 450
$ws.Range("F20").Value = 75
$ws.Range("F21").Value = 689
$ws.Range("F24").Value = 71
$ws.Range("F25").Value = 327
$ws.Range("F26").Value = 1122
$ws.Range("F27").Value = 36
$ws.Range("F28").Value = 20
$ws.Range("F29").Value = 187
$ws.Range("F30").Value = 203
$ws.Range("F32").Value = 570
$ws.Range("F34").Value = 82
$ws.Range("F36").Value = 2878
$ws.Range("F37").Value = 26
$ws.Range("F40").Value = 12713
$ws.Range("F41").Value = 1341
$ws.Range("F43").Value = 32
$ws.Range("F44").Value = 61
$ws.Range("F45").Value = 276
$ws.Range("F46").Value = 381
$ws.Range("F47").Value = 4012
$ws.Range("F49").Value = 96
